$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Header row (row 1) - new columns G, H, J, I (order matters for shared-string table)
$ws.Range("G1").Value = "Serial Number"
$ws.Range("H1").Value = "Part Number"
$ws.Range("J1").Value = "Superior Equipment"
$ws.Range("I1").Value = "Status of an object"

# Status column (I) filled first across all data rows - establishes "OPER" shared string
$ws.Range("I2").Value = "OPER"
$ws.Range("I3").Value = "OPER"
$ws.Range("I4").Value = "OPER"
$ws.Range("I5").Value = "OPER"
$ws.Range("I6").Value = "OPER"
$ws.Range("I7").Value = "OPER"
$ws.Range("I8").Value = "OPER"
$ws.Range("I9").Value = "OPER"
$ws.Range("I10").Value = "OPER"
$ws.Range("I11").Value = "OPER"

# Part Number column (H)
$ws.Range("H2").Value = "GMM-003"
$ws.Range("H5").Value = "OBIE CZ"

# Serial Number column (G)
$ws.Range("G2").Value = "ZX81"
$ws.Range("G4").Value = "ZX83"
$ws.Range("G5").Value = "ZX84"
$ws.Range("G7").Value = "ZX86"
$ws.Range("G8").Value = "ZX87"
$ws.Range("G9").Value = "ZX88"
$ws.Range("G11").Value = "ZX90"

# Superior Equipment column (J) - numeric values
$ws.Range("J2").Value = 1000100001
$ws.Range("J3").Value = 1000100001
$ws.Range("J4").Value = 1000100001
$ws.Range("J7").Value = 1000100002
$ws.Range("J8").Value = 1000100002
$ws.Range("J9").Value = 1000100002

# Update the active selection to match the diff (frozen bottom pane now centers on J15)
$ws.Range("J15").Select()
